$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1361.5
$ws.Range("I15").Value = 1361.5
$ws.Range("K15").Value = 4084.5
$ws.Range("M15").Value = -3915.5
$ws.Range("H19").Value = 963.40625
$ws.Range("I19").Value = 488.7143
$ws.Range("J19").Value = 1332.6111
$ws.Range("K19").Value = 488.7143
$ws.Range("L19").Value = 1332.6111
$ws.Range("M19").Value = -313.7143
$ws.Range("N19").Value = -1682.6111
$ws.Range("H64").Value = 8006.593
$ws.Range("I64").Value = 6941.2856
$ws.Range("J64").Value = 8379.450000000001
$ws.Range("K64").Value = 6941.2856
$ws.Range("L64").Value = 8379.450000000001
$ws.Range("M64").Value = -6693.2856
$ws.Range("N64").Value = -8875.450000000001
$ws.Range("H67").Value = 8006.593
$ws.Range("I67").Value = 6941.2856
$ws.Range("J67").Value = 8379.450000000001
$ws.Range("K67").Value = 6941.2856
$ws.Range("L67").Value = 8379.450000000001
$ws.Range("M67").Value = -6083.2856
$ws.Range("N67").Value = -10095.45
$ws.Range("H88").Value = 3810.25
$ws.Range("J88").Value = 3810.25
$ws.Range("L88").Value = 3810.25
$ws.Range("N88").Value = -4622.25
$ws.Range("H91").Value = 3810.25
$ws.Range("J91").Value = 3810.25
$ws.Range("L91").Value = 3810.25
$ws.Range("N91").Value = -6618.25
$ws.Range("H92").Value = 1182.95
$ws.Range("I92").Value = 1215.2941
$ws.Range("K92").Value = 1215.2941
$ws.Range("M92").Value = 32.70589999999993
$ws.Range("H112").Value = 22412
$ws.Range("J112").Value = 27265
$ws.Range("L112").Value = 81795
$ws.Range("N112").Value = -84011
$ws.Range("H138").Value = 2409
$ws.Range("I138").Value = 1045.7368
$ws.Range("J138").Value = 4661.3477
$ws.Range("K138").Value = 3137.2104
$ws.Range("L138").Value = 13984.0431
$ws.Range("M138").Value = 2002.7896
$ws.Range("N138").Value = -24264.0431

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 1177.8889
$ws.Range("I4").Value = 1199.875
$ws.Range("J4").Value = 1002
$ws.Range("K4").Value = 1199.875
$ws.Range("L4").Value = 1002
$ws.Range("M4").Value = -1083.875
$ws.Range("N4").Value = -1234
$ws.Range("H32").Value = 6223.2607
$ws.Range("I32").Value = 3862.182
$ws.Range("K32").Value = 3862.182
$ws.Range("M32").Value = -3575.182
$ws.Range("H43").Value = 33188
$ws.Range("J43").Value = 33188
$ws.Range("L43").Value = 33188
$ws.Range("N43").Value = -33814
$ws.Range("H45").Value = 6413860.5
$ws.Range("I45").Value = 10258688
$ws.Range("J45").Value = 5813.4443
$ws.Range("K45").Value = 10258688
$ws.Range("L45").Value = 5813.4443
$ws.Range("M45").Value = -10258311
$ws.Range("N45").Value = -6567.4443
$ws.Range("H110").Value = 2316656.5
$ws.Range("I110").Value = 2527132.8
$ws.Range("K110").Value = 2527132.8
$ws.Range("M110").Value = -2525087.8

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 7151561.5
$ws.Range("I86").Value = 9101472
$ws.Range("K86").Value = 9101472
$ws.Range("M86").Value = -9100349
$ws.Range("H89").Value = 7151561.5
$ws.Range("I89").Value = 9101472
$ws.Range("K89").Value = 45507360
$ws.Range("M89").Value = -45501744
$ws.Range("H99").Value = 5106053
$ws.Range("I99").Value = 7146529.5
$ws.Range("K99").Value = 7146529.5
$ws.Range("M99").Value = -7145031.5
$ws.Range("H107").Value = 3405988
$ws.Range("I107").Value = 5497567.5
$ws.Range("K107").Value = 5497567.5
$ws.Range("M107").Value = -5495647.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 388.10526
$ws.Range("I7").Value = 261.36365
$ws.Range("J7").Value = 562.375
$ws.Range("K7").Value = 261.36365
$ws.Range("L7").Value = 562.375
$ws.Range("M7").Value = -148.36365
$ws.Range("N7").Value = -788.375
$ws.Range("H22").Value = 254.77777
$ws.Range("J22").Value = 264
$ws.Range("L22").Value = 264
$ws.Range("N22").Value = -964
$ws.Range("H31").Value = 15791.363
$ws.Range("I31").Value = 2115.625
$ws.Range("K31").Value = 2115.625
$ws.Range("M31").Value = -1820.625
$ws.Range("H34").Value = 15791.363
$ws.Range("I34").Value = 2115.625
$ws.Range("K34").Value = 2115.625
$ws.Range("M34").Value = -1913.625
$ws.Range("H57").Value = 19950
$ws.Range("J57").Value = 20000
$ws.Range("L57").Value = 20000
$ws.Range("N57").Value = -21120
$ws.Range("H58").Value = 3847.0908
$ws.Range("I58").Value = 3196.4
$ws.Range("J58").Value = 4389.3335
$ws.Range("K58").Value = 3196.4
$ws.Range("L58").Value = 4389.3335
$ws.Range("M58").Value = -2993.4
$ws.Range("N58").Value = -4795.3335
$ws.Range("H62").Value = 2998.2
$ws.Range("I62").Value = 2997
$ws.Range("K62").Value = 2997
$ws.Range("M62").Value = -2373
$ws.Range("H65").Value = 2998.2
$ws.Range("I65").Value = 2997
$ws.Range("K65").Value = 14985
$ws.Range("M65").Value = -11865
$ws.Range("H132").Value = 44621.57
$ws.Range("I132").Value = 2408.25
$ws.Range("K132").Value = 7224.75
$ws.Range("M132").Value = -4694.75
$ws.Range("H134").Value = 2815.4783
$ws.Range("I134").Value = 1940.4375
$ws.Range("J134").Value = 4815.5713
$ws.Range("K134").Value = 5821.3125
$ws.Range("L134").Value = 14446.7139
$ws.Range("M134").Value = -3286.3125
$ws.Range("N134").Value = -19516.7139
$ws.Range("H135").Value = 118605.445
$ws.Range("J135").Value = 118605.445
$ws.Range("L135").Value = 118605.445
$ws.Range("N135").Value = -128745.445
$ws.Range("H136").Value = 3847.0908
$ws.Range("I136").Value = 3196.4
$ws.Range("J136").Value = 4389.3335
$ws.Range("K136").Value = 9589.200000000001
$ws.Range("L136").Value = 13168.0005
$ws.Range("M136").Value = -7039.200000000001
$ws.Range("N136").Value = -18268.0005
$ws.Range("H138").Value = 103419.664
$ws.Range("J138").Value = 103419.664
$ws.Range("L138").Value = 103419.664
$ws.Range("N138").Value = -113699.664

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 276.11765
$ws.Range("I23").Value = 130
$ws.Range("J23").Value = 285.25
$ws.Range("K23").Value = 390
$ws.Range("L23").Value = 855.75
$ws.Range("M23").Value = -155
$ws.Range("N23").Value = -1325.75
$ws.Range("H47").Value = 469.66666
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()
$ws.Range("H86").Value = 433
$ws.Range("I86").Value = 399
$ws.Range("J86").Value = 450
$ws.Range("K86").Value = 1197
$ws.Range("L86").Value = 1350
$ws.Range("M86").Value = -11
$ws.Range("N86").Value = -3722
$ws.Range("H89").Value = 433
$ws.Range("I89").Value = 399
$ws.Range("J89").Value = 450
$ws.Range("K89").Value = 3591
$ws.Range("L89").Value = 4050
$ws.Range("M89").Value = 2337
$ws.Range("N89").Value = -15906
$ws.Range("H96").Value = 14499.5
$ws.Range("I96").Value = 9995
$ws.Range("J96").Value = 15000
$ws.Range("K96").Value = 29985
$ws.Range("L96").Value = 45000
$ws.Range("M96").Value = -27926
$ws.Range("N96").Value = -49118

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 6756.8
$ws.Range("I2").Value = 116.454544
$ws.Range("K2").Value = 116.454544
$ws.Range("M2").Value = -3.454543999999999
$ws.Range("H70").Value = 11771700
$ws.Range("I70").Value = 16673908
$ws.Range("K70").Value = 16673908
$ws.Range("M70").Value = -16673638
$ws.Range("H73").Value = 11771700
$ws.Range("I73").Value = 16673908
$ws.Range("K73").Value = 16673908
$ws.Range("M73").Value = -16672972
$ws.Range("H97").Value = 1490020.5
$ws.Range("I97").Value = 2647170.8
$ws.Range("K97").Value = 2647170.8
$ws.Range("M97").Value = -2646674.8

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4057.0557
$ws.Range("I7").Value = 2118.889
$ws.Range("J7").Value = 5995.222
$ws.Range("K7").Value = 2118.889
$ws.Range("L7").Value = 5995.222
$ws.Range("M7").Value = -2006.889
$ws.Range("N7").Value = -6219.222
$ws.Range("H22").Value = 64833.355
$ws.Range("I22").Value = 178373.6
$ws.Range("J22").Value = 1755.4445
$ws.Range("K22").Value = 178373.6
$ws.Range("L22").Value = 1755.4445
$ws.Range("M22").Value = -178078.6
$ws.Range("N22").Value = -2345.4445
$ws.Range("H27").Value = 64833.355
$ws.Range("I27").Value = 178373.6
$ws.Range("J27").Value = 1755.4445
$ws.Range("K27").Value = 178373.6
$ws.Range("L27").Value = 1755.4445
$ws.Range("M27").Value = -178266.6
$ws.Range("N27").Value = -1969.4445
$ws.Range("H93").Value = 27780336
$ws.Range("I93").Value = 30305684
$ws.Range("J93").Value = 1500
$ws.Range("K93").Value = 30305684
$ws.Range("L93").Value = 1500
$ws.Range("M93").Value = -30304436
$ws.Range("N93").Value = -3996
$ws.Range("H126").Value = 4057.0557
$ws.Range("I126").Value = 2118.889
$ws.Range("J126").Value = 5995.222
$ws.Range("K126").Value = 6356.667
$ws.Range("L126").Value = 17985.666
$ws.Range("M126").Value = -3886.667
$ws.Range("N126").Value = -22925.666
$ws.Range("H132").Value = 6897.467
$ws.Range("I132").Value = 6791.758
$ws.Range("J132").Value = 7188.1665
$ws.Range("K132").Value = 20375.274
$ws.Range("L132").Value = 21564.4995
$ws.Range("M132").Value = -17845.274
$ws.Range("N132").Value = -26624.4995

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1398.3636
$ws.Range("I113").Value = 897.8333
$ws.Range("K113").Value = 2693.4999
$ws.Range("M113").Value = -523.4998999999998
$ws.Range("H132").Value = 50051920
$ws.Range("I132").Value = 58831916
$ws.Range("J132").Value = 298629.34
$ws.Range("K132").Value = 176495748
$ws.Range("L132").Value = 895888.02
$ws.Range("M132").Value = -176493218
$ws.Range("N132").Value = -900948.02
